# Handback status report refresh: the CI run regenerated the Xliff
# handoff/handback timestamps, so the "Latest HO Xliff Generate Date" /
# "Correspond Handoff Datetime" / "Correspond Handback DateTime" cells
# on each sheet get bumped to the new run's timestamps.

$wb = $excel.ActiveWorkbook

# Overview sheet: G2 "Latest HO Xliff Generate Date" for the first row.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-29 01:03:16"

# zh-cn sheet: H2 "Correspond Handoff Datetime", K2 "Correspond Handback DateTime".
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-29 01:03:11"
$wsZhCn.Range("K2").Value = "2016-08-29 01:03:28"

# de-de sheet: H2 "Correspond Handoff Datetime" (same value as Overview!G2),
# K2 "Correspond Handback DateTime".
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-29 01:03:16"
$wsDeDe.Range("K2").Value = "2016-08-29 01:03:35"
